$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.747.87'
$ws.Range('D2').ClearFormats()

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.740.85'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -5.17%  '
$ws.Range('E3').ClearFormats()

# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E4').ClearFormats()

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.02'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -8.70%  '
$ws.Range('E5').ClearFormats()

# Row 6
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E6').ClearFormats()

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5046'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -6.16%  '
$ws.Range('E7').ClearFormats()

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '41.83'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -6.75%  '
$ws.Range('E8').ClearFormats()

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2644'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -12.44%  '
$ws.Range('E9').ClearFormats()

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06150'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -10.50%  '
$ws.Range('E10').ClearFormats()

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.741.62'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -5.25%  '
$ws.Range('E11').ClearFormats()

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06915'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -4.36%  '
$ws.Range('E12').ClearFormats()

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.31'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -12.90%  '
$ws.Range('E13').ClearFormats()

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.490'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -9.63%  '
$ws.Range('E14').ClearFormats()

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5954'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -19.28%  '
$ws.Range('E15').ClearFormats()

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '76.70'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -14.10%  '
$ws.Range('E16').ClearFormats()

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.001'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('E17').ClearFormats()

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '25.745.74'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('E19').ClearFormats()

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000006805'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -13.65%  '
$ws.Range('E20').ClearFormats()

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.60'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -16.10%  '
$ws.Range('E21').ClearFormats()

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.964.96'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -5.64%  '
$ws.Range('E22').ClearFormats()

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.043'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -11.65%  '
$ws.Range('E23').ClearFormats()

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.172'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -13.23%  '
$ws.Range('E24').ClearFormats()

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.100'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -12.40%  '
$ws.Range('E25').ClearFormats()

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '138.05'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -3.19%  '
$ws.Range('E26').ClearFormats()

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.513'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -10.30%  '
$ws.Range('E27').ClearFormats()

# Row 28
$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('B28').ClearFormats()
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').ClearFormats()
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '14.95'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -11.81%  '
$ws.Range('E28').ClearFormats()

# Row 29
$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('B29').ClearFormats()
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C29').ClearFormats()
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.804'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -17.66%  '
$ws.Range('E29').ClearFormats()

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '103.24'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -6.49%  '
$ws.Range('E30').ClearFormats()

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.755'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -10.98%  '
$ws.Range('E31').ClearFormats()

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08089'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -8.10%  '
$ws.Range('E32').ClearFormats()

# Row 33
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -13.70%  '
$ws.Range('E33').ClearFormats()

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04498'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -6.39%  '
$ws.Range('E34').ClearFormats()

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9993'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E35').ClearFormats()

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.648'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -9.63%  '
$ws.Range('E36').ClearFormats()

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9760'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -13.65%  '
$ws.Range('E37').ClearFormats()

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.6071'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -16.75%  '
$ws.Range('E38').ClearFormats()

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.651'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -14.28%  '
$ws.Range('E39').ClearFormats()

# Row 40
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -9.38%  '
$ws.Range('E40').ClearFormats()

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9999'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E41').ClearFormats()

# Row 42
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -17.04%  '
$ws.Range('E42').ClearFormats()

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '102.85'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -4.68%  '
$ws.Range('E43').ClearFormats()

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3798'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -19.50%  '
$ws.Range('E44').ClearFormats()

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.086'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -13.46%  '
$ws.Range('E45').ClearFormats()

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.7325'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -19.28%  '
$ws.Range('E46').ClearFormats()

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05345'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -7.62%  '
$ws.Range('E47').ClearFormats()

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.1111'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -9.84%  '
$ws.Range('E48').ClearFormats()

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '30.12'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -13.26%  '
$ws.Range('E49').ClearFormats()

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.887'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -19.98%  '
$ws.Range('E50').ClearFormats()

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '52.47'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -12.43%  '
$ws.Range('E51').ClearFormats()
